$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'264.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'1.24%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'26.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-1.95%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.699"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.45%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.06110"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-1.17%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'6.735"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'0.99%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.8505"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.09%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9100"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-0.91%"
$ws.Range("E8").ClearFormats()
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.04769"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-1.86%"
$ws.Range("E9").ClearFormats()
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07099"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'0.19%"
$ws.Range("E10").ClearFormats()
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03129"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'0.57%"
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09045"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'0.03%"
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001537"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.16%"
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006168"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.26%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005973"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-0.48%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.452"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.08%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.163"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.32%"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.146"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-0.64%"
$ws.Range("E18").ClearFormats()
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3072"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-0.19%"
$ws.Range("E19").ClearFormats()
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "'0.1408"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-0.25%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-1.42%"
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'1.14%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04242"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.02%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001175"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-3.23%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004067"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'6.97%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("E26").ClearFormats()
$ws.Range("D40").Value = "'0.03936"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'1.56%"
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'0.20%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'1.83%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.002109"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-3.77%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.01148"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-29.66%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005065"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-2.04%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E46").ClearFormats()
$ws.Range("D48").Value = "'0.2577"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'58.72%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").ClearFormats()
